# edit.ps1 -- apply the commit's changes via PowerPoint COM-interop
#
# 1) Slide 6's table switches to a different table style (tableStyleId).
# 2) The deck's two embedded themes ("Office Theme" in ppt/theme/theme1.xml
#    and "Integral" in ppt/theme/theme2.xml) swap places: the colour
#    scheme that is actually in effect for the slide master/slides
#    (exposed as ppt/theme/theme2.xml) becomes the "Office" palette that
#    used to live in theme1.xml, while the notes-only theme1.xml keeps
#    its structure (its colours are not reachable through the PowerPoint
#    object model, which only ever exposes the master/slide theme).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$table  = $slide6.Shapes.Item(2).Table
$table.ApplyStyle("{28C94E99-89E6-49EE-AF5E-CD1F2AC4548E}")

# --- 2) Theme colour swap --------------------------------------------------
# Converts a "RRGGBB" hex string into the BGR-packed long that the
# PowerPoint object model's RGB properties expect.
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the "Office Theme" colours that used to be theme1.xml's
# <a:clrScheme>, now applied as the live theme (theme2.xml) colour scheme,
# in the 12-slot ThemeColorScheme order (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink).
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $slide6.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgb($officeColors[$i - 1])
}
